$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Reorder comma-separated IA Control identifiers in column A for the rows below.
$ws.Range("A2").Value = 'AU-4,AU-4 (1)'
$ws.Range("A3").Value = 'AU-14 (1),AU-4'
$ws.Range("A4").Value = 'AU-4,CM-6 b'
$ws.Range("A5").Value = 'SC-5 (2),CM-6 b,SC-5'
$ws.Range("A6").Value = 'AU-12 (3),AU-7 b,AU-8 b,AC-6 (8),AC-6 (9),CM-5 (1),AU-7 a'
$ws.Range("A7").Value = 'AU-12 (3),CM-6 b,AU-12 c,AU-7 b,AU-8 b,AU-12 a,CM-5 (1),AU-7 a'
$ws.Range("A13").Value = 'CM-7 (5) (b),CM-7 (2)'
$ws.Range("A14").Value = 'CM-7 (5) (b),CM-7 (2)'
$ws.Range("A16").Value = 'CM-7 (2),CM-6 b'
$ws.Range("A21").Value = 'CM-7 (2),CM-6 b'
$ws.Range("A22").Value = 'CM-7 (2),CM-6 b'
$ws.Range("A37").Value = 'AC-7 b,AC-7 a'
$ws.Range("A38").Value = 'AC-7 b,AC-7 a'
$ws.Range("A39").Value = 'AC-7 b,AC-7 a'
$ws.Range("A40").Value = 'AC-7 b,AC-7 a'
$ws.Range("A44").Value = 'IA-8,IA-2,AU-3 (1)'
$ws.Range("A45").Value = 'AU-3 (1),AU-12 c,MA-4 (1) (a),AU-12 a,AU-3'
$ws.Range("A46").Value = 'AU-3 (1),AU-12 c,MA-4 (1) (a),AU-12 a,AU-3'
$ws.Range("A47").Value = 'AU-3 (1),AU-12 c,MA-4 (1) (a),AU-12 a,AU-3'
$ws.Range("A48").Value = 'AU-3 (1),AU-12 c,MA-4 (1) (a),AU-12 a,AU-3'
$ws.Range("A49").Value = 'AU-3 (1),AU-12 c,MA-4 (1) (a),AU-12 a,AU-3'
$ws.Range("A50").Value = 'AU-3 (1),AU-12 c,MA-4 (1) (a),AU-12 a,AU-3'
$ws.Range("A51").Value = 'AU-3 (1),AU-12 c,MA-4 (1) (a),AU-12 a,AU-3'
$ws.Range("A52").Value = 'AU-3 (1),AU-12 c,MA-4 (1) (a),AU-12 a,AU-3'
$ws.Range("A53").Value = 'AU-3 (1),AU-12 c,MA-4 (1) (a),AU-12 a,AU-3'
$ws.Range("A54").Value = 'AU-3 (1),AU-12 c,MA-4 (1) (a),AU-12 a,AU-3'
$ws.Range("A55").Value = 'AU-3 (1),AU-12 c,MA-4 (1) (a),AU-12 a,AU-3'
$ws.Range("A56").Value = 'AU-3 (1),AU-12 c,MA-4 (1) (a),AU-12 a,AU-3'
$ws.Range("A57").Value = 'AU-3 (1),AU-12 c,MA-4 (1) (a),AU-12 a,AU-3'
$ws.Range("A58").Value = 'AU-3 (1),AU-12 c,MA-4 (1) (a),AU-12 a,AU-3'
$ws.Range("A59").Value = 'AU-3 (1),AU-12 c,MA-4 (1) (a),AU-12 a,AU-3'
$ws.Range("A60").Value = 'AU-3 (1),AU-12 c,MA-4 (1) (a),AU-12 a,AU-3'
$ws.Range("A61").Value = 'AU-3 (1),AU-12 c,MA-4 (1) (a),AU-12 a,AU-3'
$ws.Range("A62").Value = 'AU-3 (1),AU-12 c,MA-4 (1) (a),AU-12 a,AU-3'
$ws.Range("A63").Value = 'AU-3 (1),AU-12 c,MA-4 (1) (a),AU-12 a,AU-3'
$ws.Range("A64").Value = 'AU-3 (1),AU-12 c,MA-4 (1) (a),AU-12 a,AU-3'
$ws.Range("A65").Value = 'AU-3 (1),AU-12 c,MA-4 (1) (a),AU-12 a,AU-3'
$ws.Range("A66").Value = 'AU-3 (1),AU-12 c,MA-4 (1) (a),AU-12 a,AU-3'
$ws.Range("A67").Value = 'AU-3 (1),AU-12 c,MA-4 (1) (a),AU-12 a,AU-3'
$ws.Range("A68").Value = 'AU-3 (1),AU-12 c,MA-4 (1) (a),AU-12 a,AU-3'
$ws.Range("A69").Value = 'AU-3 (1),AU-12 c,MA-4 (1) (a),AU-12 a,AU-3'
$ws.Range("A70").Value = 'AU-3 (1),AU-12 c,MA-4 (1) (a),AU-12 a,AU-3'
$ws.Range("A71").Value = 'AU-3 (1),AU-12 c,MA-4 (1) (a),AU-12 a,AU-3'
$ws.Range("A72").Value = 'AU-3 (1),AU-12 c,MA-4 (1) (a),AU-12 a,AU-3'
$ws.Range("A73").Value = 'AU-3 (1),AU-12 c,MA-4 (1) (a),AU-12 a,AU-3'
$ws.Range("A74").Value = 'AU-3 (1),AU-12 c,MA-4 (1) (a),AU-12 a,AU-3'
$ws.Range("A75").Value = 'AU-3 (1),AU-12 c,MA-4 (1) (a),AU-12 a,AU-3'
$ws.Range("A76").Value = 'AU-3 (1),AU-12 c,MA-4 (1) (a),AU-12 a,AU-3'
$ws.Range("A77").Value = 'AU-3 (1),AU-12 c,MA-4 (1) (a),AU-12 a,AU-3'
$ws.Range("A78").Value = 'AU-3 (1),AU-12 c,MA-4 (1) (a),AU-12 a,AU-3'
$ws.Range("A79").Value = 'AU-3 (1),AU-12 c,MA-4 (1) (a),AU-12 a,AU-3'
$ws.Range("A80").Value = 'AU-3 (1),AU-12 c,MA-4 (1) (a),AU-12 a,AU-3'
$ws.Range("A81").Value = 'AU-3 (1),AU-12 c,MA-4 (1) (a),AU-12 a,AU-3'
$ws.Range("A82").Value = 'AU-3 (1),AU-12 c,MA-4 (1) (a),AU-12 a,AU-3'
$ws.Range("A83").Value = 'AU-3 (1),AU-12 c,MA-4 (1) (a),AU-12 a,AU-3'
$ws.Range("A84").Value = 'AU-3,AU-3 (1),MA-4 (1) (a)'
$ws.Range("A85").Value = 'AU-3 (1),AU-12 c,MA-4 (1) (a),AU-12 a,AU-3'
$ws.Range("A86").Value = 'AU-3 (1),AU-12 c,MA-4 (1) (a),AU-12 a,AU-3'
$ws.Range("A87").Value = 'AU-3 (1),AU-12 c,MA-4 (1) (a),AU-12 a,AU-3'
$ws.Range("A88").Value = 'AU-3 (1),AU-12 c,MA-4 (1) (a),AU-12 a,AU-3'
$ws.Range("A89").Value = 'AU-3 (1),AU-12 c,MA-4 (1) (a),AU-12 a,AU-3'
$ws.Range("A91").Value = 'AU-3 (1),AU-12 c,MA-4 (1) (a),AU-12 a,AU-3'
$ws.Range("A92").Value = 'AU-3 (1),AU-12 c,MA-4 (1) (a),AU-12 a,AU-3'
$ws.Range("A93").Value = 'AU-3 (1),AU-12 c,MA-4 (1) (a),AU-12 a,AU-3'
$ws.Range("A94").Value = 'AU-3 (1),AU-12 c,MA-4 (1) (a),AU-12 a,AU-3'
$ws.Range("A95").Value = 'AU-3 (1),AU-12 c,MA-4 (1) (a),AU-12 a,AU-3'
$ws.Range("A96").Value = 'AU-12 c,AU-3,AU-3 (1),MA-4 (1) (a)'
$ws.Range("A97").Value = 'AU-3 (1),AU-12 c,MA-4 (1) (a),AU-12 a,AU-3'
$ws.Range("A98").Value = 'AU-3 (1),AU-12 c,MA-4 (1) (a),AU-12 a,AU-3'
$ws.Range("A99").Value = 'AU-3 (1),AU-12 c,MA-4 (1) (a),AU-12 a,AU-3'
$ws.Range("A100").Value = 'AU-3 (1),AU-12 c,MA-4 (1) (a),AU-12 a,AU-3'
$ws.Range("A101").Value = 'AU-3 (1),AU-12 c,MA-4 (1) (a),AU-12 a,AU-3'
$ws.Range("A102").Value = 'AC-2 (4),AU-3 (1),AU-12 c,MA-4 (1) (a),AU-12 a,AU-3'
$ws.Range("A103").Value = 'AC-2 (4),AU-3 (1),AU-12 c,MA-4 (1) (a),AU-12 a,AU-3'
$ws.Range("A104").Value = 'AC-2 (4),AU-3 (1),AU-12 c,MA-4 (1) (a),AU-12 a,AU-3'
$ws.Range("A105").Value = 'AC-2 (4),AU-3 (1),AU-12 c,MA-4 (1) (a),AU-3'
$ws.Range("A106").Value = 'AC-2 (4),AU-3 (1),AU-12 c,MA-4 (1) (a),AU-12 a,AU-3'
$ws.Range("A107").Value = 'AC-2 (4),AU-3 (1),AU-12 c,MA-4 (1) (a),AU-12 a,AU-3'
$ws.Range("A108").Value = 'AC-2 (4),AU-3 (1),AU-12 c,MA-4 (1) (a),AU-12 a,AU-3'
$ws.Range("A109").Value = 'AC-2 (4),AU-3 (1),AU-12 c,MA-4 (1) (a),AU-12 a,AU-3'
$ws.Range("A110").Value = 'AC-2 (4),AU-3 (1),AU-12 c,MA-4 (1) (a),AU-12 a,AU-3'
$ws.Range("A111").Value = 'AU-3 (1),AU-12 c,AU-14 (1),MA-4 (1) (a),AU-12 a,AU-3'
$ws.Range("A119").Value = 'AU-12 a,AU-12 c,AU-3,MA-4 (1) (a)'
$ws.Range("A120").Value = 'AU-12 a,AU-12 c,AU-3,MA-4 (1) (a)'
$ws.Range("A121").Value = 'AU-9,AU-12 c'
$ws.Range("A126").Value = 'AC-2 (4),CM-5 (1),AC-6 (9),AU-12 c'
$ws.Range("A128").Value = 'IA-5 (1) (a),IA-5 (1) (b),CM-6 b'
$ws.Range("A132").Value = 'SC-8,SC-13,AC-17 (2),MA-4 c'
$ws.Range("A133").Value = 'MA-4 e,MA-4 (7),SC-10,AC-12'
$ws.Range("A136").Value = 'AC-11 a,SC-10'
$ws.Range("A137").Value = 'AU-3 (1),CM-6 b,AU-14 (1),MA-4 (1) (a),AU-12 a,AU-6 (4),CM-5 (1),AU-3,AU-7 (1),AU-7 a'
$ws.Range("A142").Value = 'MA-4 (1) (a),AU-12 c'
$ws.Range("A143").Value = 'MA-4 (1) (a),AU-12 c'
$ws.Range("A144").Value = 'MA-4 (1) (a),AU-12 c'
$ws.Range("A145").Value = 'MA-4 (1) (a),AU-12 c'
$ws.Range("A146").Value = 'MA-4 (1) (a),AU-12 c'
$ws.Range("A147").Value = 'MA-4 (1) (a),AU-12 c'
$ws.Range("A148").Value = 'MA-4 (1) (a),AU-12 c'
$ws.Range("A149").Value = 'MA-4 (1) (a),AU-12 c'
$ws.Range("A150").Value = 'MA-4 (1) (a),AU-12 c'
$ws.Range("A151").Value = 'MA-4 (1) (a),AU-12 c'
$ws.Range("A152").Value = 'MA-4 (1) (a),AU-12 c'
$ws.Range("A153").Value = 'MA-4 (1) (a),AU-12 c'
$ws.Range("A154").Value = 'MA-4 (1) (a),AU-12 c'
$ws.Range("A165").Value = 'SC-8 (1),SC-8,SC-8 (2)'
$ws.Range("A166").Value = 'SC-8 (1),SC-8,SC-8 (2)'
$ws.Range("A167").Value = 'AC-17 (2),SC-8'
$ws.Range("A169").Value = 'AC-17 (2),SC-8'
$ws.Range("A171").Value = 'AC-11 a,AC-11 b'
$ws.Range("A178").Value = 'CM-6 b,AU-4 (1),AU-6 (4)'
$ws.Range("A180").Value = 'AC-17 (1),CM-6 b,CM-7 b'
$ws.Range("A181").Value = 'SI-11 b,AU-9'
$ws.Range("A182").Value = 'SI-11 b,AU-9'
$ws.Range("A183").Value = 'SI-11 b,AU-9'
$ws.Range("A184").Value = 'SI-11 b,AU-9'
$ws.Range("A185").Value = 'SI-11 b,AU-9'
$ws.Range("A192").Value = 'AU-3,CM-6 b'
$ws.Range("A198").Value = 'AU-3,AU-4 (1)'
$ws.Range("A214").Value = 'AC-2 (4),AC-6 (9),AU-12 c'
$ws.Range("A220").Value = 'IA-2,IA-2 (5),IA-2 (3),IA-2 (2),IA-2 (4)'
$ws.Range("A221").Value = 'IA-2,IA-2 (5),IA-2 (3),IA-2 (2),IA-2 (4)'
$ws.Range("A222").Value = 'AC-18 (1),SC-8 (1),SC-8'
$ws.Range("A225").Value = 'IA-7,CM-6 b'
$ws.Range("A226").Value = 'IA-7,CM-6 b'
$ws.Range("A227").Value = 'IA-7,CM-6 b'
$ws.Range("A243").Value = 'CM-6 b,SC-2,SI-16'
$ws.Range("A257").Value = 'IA-3,CM-6 b'
$ws.Range("A258").Value = 'IA-3,CM-6 b'
$ws.Range("A259").Value = 'IA-3,CM-6 b'
$ws.Range("A260").Value = 'IA-3,CM-6 b'
$ws.Range("A269").Value = 'IA-2 (1),IA-2 (3),IA-2 (2),IA-2 (4)'
$ws.Range("A274").Value = 'SC-4,CM-6 b'
$ws.Range("A279").Value = 'AU-12 a,CM-6 b'
$ws.Range("A299").Value = 'IA-2 (1),IA-2 (11),IA-2 (12)'
$ws.Range("A309").Value = 'AU-8 b,AU-8 (1) (b),AU-8 (1) (a)'
$ws.Range("A328").Value = 'CM-5 (1),AU-12 c'
$ws.Range("A342").Value = 'IA-3,CM-7 b'
$ws.Range("A347").Value = 'CM-6 b,IA-5 (1) (c),CM-7 a'
$ws.Range("A358").Value = 'AC-11 (1),AC-11 b'
$ws.Range("A361").Value = 'CM-3 (5),SI-6 d,SI-6 b'
$ws.Range("A362").Value = 'CM-7 a,CM-6 b'
$ws.Range("A385").Value = 'AC-17 (2),CM-6 b'
$ws.Range("A389").Value = 'SC-3,SI-6 a'
$ws.Range("A391").Value = 'IA-5 (1) (a),CM-6 b'
$ws.Range("A401").Value = 'SC-3,CM-6 b'
$ws.Range("A402").Value = 'SC-3,CM-6 b'
$ws.Range("A403").Value = 'SC-3,CM-6 b'
$ws.Range("A448").Value = 'IA-5 (1) (c),CM-6 b'
$ws.Range("A541").Value = 'SI-2 (2),CM-6 b'
$ws.Range("A550").Value = 'SI-2 (2),CM-6 b'
$ws.Range("A558").Value = 'CM-3 (5),SI-6 a'

# Populate the Fix (M) cell for row 131, which was previously empty.
$m131 = 'Configure Red Hat Enterprise Linux 9 to run in FIPS mode.
Run the following commands:
$ sudo fips-mode-setup --enable
$ sudo update-crypto-policies --set 
The system needs to be rebooted for these changes to take effect.'
$ws.Range("M131").Value = $m131

